$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right under the title, and re-insert its content as a new
#    paragraph right before the DALLE image-prompt paragraph at the
#    end of the document (Cut + Paste preserves the run / formatting
#    structure, incl. the bold "Meta description" run).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$metaRange.Cut()

$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
$pasteTarget = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$pasteTarget.Paste()

# ------------------------------------------------------------------
# 2) The pasted paragraph (now second-to-last) still reads:
#      "Meta description: Read our review of Beat the Beast Mighty
#       Sphinx and play for free. Enjoy simple gameplay, high
#       volatility, and impressive graphics."
#    Turn the bold "Meta description" run into the new bold headline
#    "Play Beat the Beast Mighty Sphinx Free - Review", then strip the
#    leftover ": Read our review..." tail from that paragraph (it
#    becomes the new text of the image-prompt paragraph below).
# ------------------------------------------------------------------
$newHeadlinePara = $d.Paragraphs.Item($paraCount)
$headlineRange = $d.Range($newHeadlinePara.Range.Start, $newHeadlinePara.Range.End)
$headlineRange.Find.Execute("Meta description", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Play Beat the Beast Mighty Sphinx Free - Review", 2)

$newHeadlinePara2 = $d.Paragraphs.Item($paraCount)
$tailRange = $d.Range($newHeadlinePara2.Range.Start, $newHeadlinePara2.Range.End)
$tailRange.Find.Execute(": Read our review of Beat the Beast Mighty Sphinx and play for free. Enjoy simple gameplay, high volatility, and impressive graphics.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ------------------------------------------------------------------
# 3) Replace the DALLE image-prompt text (last paragraph) with the
#    "Read our review..." sentence, keeping its italic formatting.
# ------------------------------------------------------------------
$promptPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$promptRange = $d.Range($promptPara.Range.Start, $promptPara.Range.End)
$oldPrompt = "Prompt: DALLE, create a cartoon-style feature image for " + [char]34 + "Beat the Beast Mighty Sphinx" + [char]34 + " featuring a happy Maya warrior with glasses. Description: The feature image should be in cartoon-style with bright and vibrant colors. The main element of the image should be a happy and confident Maya warrior with glasses, standing in front of a giant Sphinx. The warrior should be wearing traditional Maya clothing, with a feather headdress and accessories. The background should have an Egyptian theme, with hieroglyphics and pyramids visible. The Sphinx should be portrayed as dark and imposing, with glowing yellow eyes. The image should emphasize the adventure, excitement, and mystery of Ancient Egypt, while also showcasing the unique blend of Maya and Egyptian elements in the game."
$newPrompt = "Read our review of Beat the Beast Mighty Sphinx and play for free. Enjoy simple gameplay, high volatility, and impressive graphics."
$promptRange.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2)
